# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly generated site output.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 33
$wsExhibit.Range("F4").Value = 64
$wsExhibit.Range("F5").Value = 5088
$wsExhibit.Range("F7").Value = 32
$wsExhibit.Range("F8").Value = 90
$wsExhibit.Range("F9").Value = 314
$wsExhibit.Range("F10").Value = 57

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 33
$wsAll.Range("F8").Value = 64
$wsAll.Range("F9").Value = 5088
$wsAll.Range("F11").Value = 32
$wsAll.Range("F12").Value = 90
$wsAll.Range("F14").Value = 314
$wsAll.Range("F15").Value = 57
